# Saldo.xlsx update script
# 1. Update row 3 (account 005395948 / THAIS / 100000) -> (004813166 / VENIA / 123344.17)
# 2. Delete row containing account 004399832 / EULER / 16143.27
# 3. Delete two rows containing accounts 004813166 / VENIA / 10782.39 and 005654767 / DIEGO / 9100

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 in place with the new account/name/balance.
# The account number must stay text with its leading zeros intact. Typing
# it straight into Value would auto-convert to a number, so stage it in a
# scratch cell (forced to text via a leading apostrophe), copy only the
# value across with PasteSpecial, then remove the scratch cell again.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Value = "'004813166"
$scratch.Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item(3, 2).Value = "VENIA"
$ws.Cells.Item(3, 3).Value = 123344.17

# Delete the row for account 004399832 / EULER (originally row 13)
$ws.Rows.Item(13).Delete()

# After deleting row 13, the rows that were 20 and 21 shift up to 19 and 20
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()
